$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 64

# Copy formatting (styles) from the previous data row so the new row
# matches the existing look & feel (bold/centered index column, date format, etc.)
$ws.Range("A63:V63").Copy($ws.Range("A64:V64"))

$ws.Cells.Item($row, 1).Value = 63
$ws.Cells.Item($row, 2).Value = "azerbaijan"
$ws.Cells.Item($row, 3).Value = "premier-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45242.5625
$ws.Cells.Item($row, 6).Value = "Zira"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Gabala"
$ws.Cells.Item($row, 9).Value = 0

$ws.Cells.Item($row, 10).Value = 2.08
$ws.Cells.Item($row, 11).Value = "11/11/2023 01:42"
$ws.Cells.Item($row, 12).Value = 1.92
$ws.Cells.Item($row, 13).Value = "12/11/2023 13:22"

$ws.Cells.Item($row, 14).Value = 2.91
$ws.Cells.Item($row, 15).Value = "11/11/2023 01:42"
$ws.Cells.Item($row, 16).Value = 3.08
$ws.Cells.Item($row, 17).Value = "12/11/2023 13:22"

$ws.Cells.Item($row, 18).Value = 3.45
$ws.Cells.Item($row, 19).Value = "11/11/2023 01:42"
$ws.Cells.Item($row, 20).Value = 4.43
$ws.Cells.Item($row, 21).Value = "12/11/2023 13:22"

$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/zira-fk-gabala/2HhQQozN/"
